$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Update the shared-string values that changed in the "Nome"/"Sobrenome" columns
$ws.Range("A2").Value = "Itqachi"
$ws.Range("A3").Value = "sasuke"
$ws.Range("A4").Value = "oii"
$ws.Range("A5").Value = "Oláaa"
$ws.Range("A7").Value = "Testeee"
$ws.Range("A8").Value = "Gabriel"
$ws.Range("B8").Value = "Testando"

# Reset the stale selection (was H20) back to the default top-left cell
$ws.Range("A1").Select()
